$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("A200")
$cell.Value = "Test1"
$cell.Interior.Color = 65535
$cell.Font.Bold = $true
$cell.Borders.Item(7).LineStyle = 1
$cell.Borders.Item(10).LineStyle = 1
$cell.Borders.Item(8).LineStyle = 1

$cell2 = $ws.Range("C300")
$cell2.Value = "Test2"
$cell2.Interior.Color = 65535
$cell2.Font.Bold = $true
$cell2.Borders.Item(7).LineStyle = 1
$cell2.Borders.Item(10).LineStyle = 1
$cell2.Borders.Item(8).LineStyle = 1

Write-Host "Cell1 L,T,B,R:" $cell.Borders.Item(7).LineStyle $cell.Borders.Item(8).LineStyle $cell.Borders.Item(9).LineStyle $cell.Borders.Item(10).LineStyle
Write-Host "Cell2 L,T,B,R:" $cell2.Borders.Item(7).LineStyle $cell2.Borders.Item(8).LineStyle $cell2.Borders.Item(9).LineStyle $cell2.Borders.Item(10).LineStyle
